# Automatic update from scheduled task (Actualización automática desde tarea programada)
# Corrects the timestamp of row 11 and appends a new reading as row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: fix the slightly-off execution timestamp in column A
$ws.Range("A11").Value2 = 45874.37516971065

# Row 12: new sensor reading appended by the scheduled task
$ws.Range("A12").Value2 = 45874.41689232497
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat

$ws.Range("B12").Value2 = 2025
$ws.Range("C12").Value2 = 19
$ws.Range("D12").Value2 = 14.78
$ws.Range("E12").Value2 = 92.93000000000001
$ws.Range("F12").Value2 = 186.33
$ws.Range("G12").Value2 = 7.28
$ws.Range("H12").Value2 = "ESE"
$ws.Range("I12").Value2 = 0
$ws.Range("J12").Value2 = "10:00:19"
